$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H53").Value = 1464.1333
$ws_ALC.Range("I53").Value = 1047
$ws_ALC.Range("K53").Value = 1047
$ws_ALC.Range("M53").Value = -410
$ws_ALC.Range("H54").Value = 4000
$ws_ALC.Range("I54").Value = 4000
$ws_ALC.Range("K54").Value = 4000
$ws_ALC.Range("M54").Value = -3514
$ws_ALC.Range("H129").Value = 879.4050999999999
$ws_ALC.Range("J129").Value = 900.24
$ws_ALC.Range("L129").Value = 2700.72
$ws_ALC.Range("N129").Value = -12700.72
$ws_ALC.Range("H137").Value = 1185.8572
$ws_ALC.Range("I137").Value = 1155.875
$ws_ALC.Range("J137").Value = 1281.8
$ws_ALC.Range("K137").Value = 3467.625
$ws_ALC.Range("L137").Value = 3845.4
$ws_ALC.Range("M137").Value = -917.625
$ws_ALC.Range("N137").Value = -8945.4
$ws_ALC.Range("H138").Value = 1331.88
$ws_ALC.Range("I138").Value = 677.44446
$ws_ALC.Range("J138").Value = 1700
$ws_ALC.Range("K138").Value = 2032.33338
$ws_ALC.Range("L138").Value = 5100
$ws_ALC.Range("M138").Value = 3107.66662
$ws_ALC.Range("N138").Value = -15380
$ws_ALC.Range("H141").Value = 865
$ws_ALC.Range("I141").Value = 865
$ws_ALC.Range("K141").Value = 2595
$ws_ALC.Range("M141").Value = 2585

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H60").Value = 12333.333
$ws_ARM.Range("I60").Value = 5000
$ws_ARM.Range("J60").Value = 16000
$ws_ARM.Range("K60").Value = 5000
$ws_ARM.Range("L60").Value = 16000
$ws_ARM.Range("M60").Value = -4267
$ws_ARM.Range("N60").Value = -17466
$ws_ARM.Range("H74").Value = 1015.5714
$ws_ARM.Range("I74").Value = 862.9231
$ws_ARM.Range("J74").Value = 3000
$ws_ARM.Range("K74").Value = 862.9231
$ws_ARM.Range("L74").Value = 3000
$ws_ARM.Range("M74").Value = 11.07690000000002
$ws_ARM.Range("N74").Value = -4748
$ws_ARM.Range("H76").Value = 22400
$ws_ARM.Range("J76").Value = 22400
$ws_ARM.Range("L76").Value = 22400
$ws_ARM.Range("N76").Value = -23076
$ws_ARM.Range("H77").Value = 1015.5714
$ws_ARM.Range("I77").Value = 862.9231
$ws_ARM.Range("J77").Value = 3000
$ws_ARM.Range("K77").Value = 4314.6155
$ws_ARM.Range("L77").Value = 15000
$ws_ARM.Range("M77").Value = 53.38450000000012
$ws_ARM.Range("N77").Value = -23736
$ws_ARM.Range("H79").Value = 22400
$ws_ARM.Range("J79").Value = 22400
$ws_ARM.Range("L79").Value = 22400
$ws_ARM.Range("N79").Value = -24740
$ws_ARM.Range("H96").Value = 22000
$ws_ARM.Range("J96").Value = 22000
$ws_ARM.Range("L96").Value = 22000
$ws_ARM.Range("N96").Value = -27492
$ws_ARM.Range("H102").Value = 8334237.5
$ws_ARM.Range("I102").Value = 8772845
$ws_ARM.Range("K102").Value = 8772845
$ws_ARM.Range("M102").Value = -8771223
$ws_ARM.Range("H135").Value = 17174.777
$ws_ARM.Range("J135").Value = 17174.777
$ws_ARM.Range("L135").Value = 17174.777
$ws_ARM.Range("N135").Value = -27314.777

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H100").Value = 0
$ws_BSM.Range("J100").Value = 0
$ws_BSM.Range("L100").Value = 0
$ws_BSM.Range("N100").ClearContents()
$ws_BSM.Range("H105").Value = 166669550
$ws_BSM.Range("I105").Value = 250003070
$ws_BSM.Range("J105").Value = 2511
$ws_BSM.Range("K105").Value = 250003070
$ws_BSM.Range("L105").Value = 2511
$ws_BSM.Range("M105").Value = -250001323
$ws_BSM.Range("N105").Value = -6005

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2175.76
$ws_CRP.Range("I31").Value = 1127.091
$ws_CRP.Range("J31").Value = 2999.7144
$ws_CRP.Range("K31").Value = 1127.091
$ws_CRP.Range("L31").Value = 2999.7144
$ws_CRP.Range("M31").Value = -832.0909999999999
$ws_CRP.Range("N31").Value = -3589.7144
$ws_CRP.Range("H34").Value = 2175.76
$ws_CRP.Range("I34").Value = 1127.091
$ws_CRP.Range("J34").Value = 2999.7144
$ws_CRP.Range("K34").Value = 1127.091
$ws_CRP.Range("L34").Value = 2999.7144
$ws_CRP.Range("M34").Value = -925.0909999999999
$ws_CRP.Range("N34").Value = -3403.7144
$ws_CRP.Range("H74").Value = 24820.75
$ws_CRP.Range("I74").Value = 3285
$ws_CRP.Range("J74").Value = 31999.334
$ws_CRP.Range("K74").Value = 3285
$ws_CRP.Range("L74").Value = 31999.334
$ws_CRP.Range("M74").Value = -2411
$ws_CRP.Range("N74").Value = -33747.334
$ws_CRP.Range("H77").Value = 24820.75
$ws_CRP.Range("I77").Value = 3285
$ws_CRP.Range("J77").Value = 31999.334
$ws_CRP.Range("K77").Value = 9855
$ws_CRP.Range("L77").Value = 95998.00199999999
$ws_CRP.Range("M77").Value = -5487
$ws_CRP.Range("N77").Value = -104734.002
$ws_CRP.Range("H80").Value = 14690
$ws_CRP.Range("I80").Value = 14690
$ws_CRP.Range("J80").Value = 0
$ws_CRP.Range("K80").Value = 14690
$ws_CRP.Range("L80").Value = 0
$ws_CRP.Range("M80").Value = -13567
$ws_CRP.Range("N80").ClearContents()
$ws_CRP.Range("H83").Value = 14690
$ws_CRP.Range("I83").Value = 14690
$ws_CRP.Range("J83").Value = 0
$ws_CRP.Range("K83").Value = 44070
$ws_CRP.Range("L83").Value = 0
$ws_CRP.Range("M83").Value = -38454
$ws_CRP.Range("N83").ClearContents()
$ws_CRP.Range("H116").Value = 0
$ws_CRP.Range("J116").Value = 0
$ws_CRP.Range("L116").Value = 0
$ws_CRP.Range("N116").ClearContents()

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H2").Value = 60
$ws_CUL.Range("I2").Value = 46.666668
$ws_CUL.Range("J2").Value = 80
$ws_CUL.Range("K2").Value = 280.000008
$ws_CUL.Range("L2").Value = 480
$ws_CUL.Range("M2").Value = -167.000008
$ws_CUL.Range("N2").Value = -706
$ws_CUL.Range("H34").Value = 1211.1177
$ws_CUL.Range("I34").Value = 511.125
$ws_CUL.Range("J34").Value = 1833.3334
$ws_CUL.Range("K34").Value = 1533.375
$ws_CUL.Range("L34").Value = 5500.0002
$ws_CUL.Range("M34").Value = -1449.375
$ws_CUL.Range("N34").Value = -5668.0002
$ws_CUL.Range("H39").Value = 1659.8235
$ws_CUL.Range("J39").Value = 1738.375
$ws_CUL.Range("L39").Value = 5215.125
$ws_CUL.Range("N39").Value = -5803.125
$ws_CUL.Range("H55").Value = 3500
$ws_CUL.Range("J55").Value = 3500
$ws_CUL.Range("L55").Value = 10500
$ws_CUL.Range("N55").Value = -10854
$ws_CUL.Range("H98").Value = 52
$ws_CUL.Range("I98").Value = 58
$ws_CUL.Range("J98").Value = 40
$ws_CUL.Range("K98").Value = 174
$ws_CUL.Range("L98").Value = 120
$ws_CUL.Range("M98").Value = 1324
$ws_CUL.Range("N98").Value = -3116
$ws_CUL.Range("H114").Value = 486.625
$ws_CUL.Range("I114").Value = 413.2857
$ws_CUL.Range("J114").Value = 1000
$ws_CUL.Range("K114").Value = 1239.8571
$ws_CUL.Range("L114").Value = 3000
$ws_CUL.Range("M114").Value = 2014.1429
$ws_CUL.Range("N114").Value = -9508
$ws_CUL.Range("H122").Value = 389.77777
$ws_CUL.Range("I122").Value = 374.5
$ws_CUL.Range("J122").Value = 402
$ws_CUL.Range("K122").Value = 3370.5
$ws_CUL.Range("L122").Value = 3618
$ws_CUL.Range("M122").Value = -920.5
$ws_CUL.Range("N122").Value = -8518

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H92").Value = 25403.75
$ws_GSM.Range("J92").Value = 25403.75
$ws_GSM.Range("L92").Value = 25403.75
$ws_GSM.Range("N92").Value = -29147.75

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 1455.2
$ws_LTW.Range("I22").Value = 1492.4375
$ws_LTW.Range("J22").Value = 1389
$ws_LTW.Range("K22").Value = 1492.4375
$ws_LTW.Range("L22").Value = 1389
$ws_LTW.Range("M22").Value = -1197.4375
$ws_LTW.Range("N22").Value = -1979
$ws_LTW.Range("H27").Value = 1455.2
$ws_LTW.Range("I27").Value = 1492.4375
$ws_LTW.Range("J27").Value = 1389
$ws_LTW.Range("K27").Value = 1492.4375
$ws_LTW.Range("L27").Value = 1389
$ws_LTW.Range("M27").Value = -1385.4375
$ws_LTW.Range("N27").Value = -1603
$ws_LTW.Range("H136").Value = 1985.5
$ws_LTW.Range("I136").Value = 1840.5714
$ws_LTW.Range("J136").Value = 3000
$ws_LTW.Range("K136").Value = 5521.7142
$ws_LTW.Range("L136").Value = 9000
$ws_LTW.Range("M136").Value = -2971.7142
$ws_LTW.Range("N136").Value = -14100

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H92").Value = 3275
$ws_WVR.Range("J92").Value = 3275
$ws_WVR.Range("L92").Value = 3275
$ws_WVR.Range("N92").Value = -8267
$ws_WVR.Range("H101").Value = 12257.286
$ws_WVR.Range("J101").Value = 12257.286
$ws_WVR.Range("L101").Value = 12257.286
$ws_WVR.Range("N101").Value = -18747.286

Write-Output "Applied 204 edits"